# Generate Report for Handback
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview sheet and on each language sheet (zh-cn, de-de).
# - Each language sheet gets its "Latest Target File" / "Latest Handback File"
#   columns (E/F) filled in with hyperlinks to the handed-back files, and the
#   "Latest Handback DateTime" column (G) is stamped with the handback time.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$mdFile  = "42de5bca-bee6-4481-acdf-58d06ea59353.md"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/c6b0a6ed0cb2e041933e5cf4d135b350156adbad/e2e/42de5bca-bee6-4481-acdf-58d06ea59353.md"

# --- Overview sheet: flip status for both rows / both language columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $statusNew
$wsOverview.Range("C2").Value = $statusNew
$wsOverview.Range("B3").Value = $statusNew
$wsOverview.Range("C3").Value = $statusNew

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhXlfFile = "42de5bca-bee6-4481-acdf-58d06ea59353.2876460dd14dd29860c1a9c1343044e64ae3d965.zh-cn.xlf"
$zhXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/35f7c27d1c9efef1a4be25bd84a4f30baa6dbd42/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/42de5bca-bee6-4481-acdf-58d06ea59353.2876460dd14dd29860c1a9c1343044e64ae3d965.zh-cn.xlf"

$wsZh.Range("B2").Value = $statusNew
$wsZh.Range("B3").Value = $statusNew

$wsZh.Hyperlinks.Add($wsZh.Range("E2"), $mdUrl, $null, $null, $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("F2"), $zhXlfUrl, $null, $null, $zhXlfFile)
$wsZh.Hyperlinks.Add($wsZh.Range("E3"), $mdUrl, $null, $null, $mdFile)
$wsZh.Hyperlinks.Add($wsZh.Range("F3"), $zhXlfUrl, $null, $null, $zhXlfFile)

$wsZh.Range("G2").Value = "2016-03-02 15:25:55"
$wsZh.Range("G3").Value = "2016-03-02 15:25:55"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$deXlfFile = "42de5bca-bee6-4481-acdf-58d06ea59353.2876460dd14dd29860c1a9c1343044e64ae3d965.de-de.xlf"
$deXlfUrl  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e19a44674988db32eee10be474b3665cb9c512a9/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/42de5bca-bee6-4481-acdf-58d06ea59353.2876460dd14dd29860c1a9c1343044e64ae3d965.de-de.xlf"

$wsDe.Range("B2").Value = $statusNew
$wsDe.Range("B3").Value = $statusNew

$wsDe.Hyperlinks.Add($wsDe.Range("E2"), $mdUrl, $null, $null, $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("F2"), $deXlfUrl, $null, $null, $deXlfFile)
$wsDe.Hyperlinks.Add($wsDe.Range("E3"), $mdUrl, $null, $null, $mdFile)
$wsDe.Hyperlinks.Add($wsDe.Range("F3"), $deXlfUrl, $null, $null, $deXlfFile)

$wsDe.Range("G2").Value = "2016-03-02 15:26:13"
$wsDe.Range("G3").Value = "2016-03-02 15:26:13"
